$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New data rows (40-44): additional William H Harsha Lake 2017/2016 G-res
# survey rows plus trailing blank/NA rows that were appended to the sheet.
# ---------------------------------------------------------------------------

# Row 40: William H Harsha Lake 2017
$ws.Range("A40").Value = "William H Harsha Lake 2017"
$ws.Range("B40").Value = "USEPA"
$ws.Range("C40").Value = 8.3194439599999992
$ws.Range("D40").Value = 34.451219510000001
$ws.Range("E40").Value = 0.095912589
$ws.Range("F40").Value = 11.629
$ws.Range("G40").Value = 0.1552
$ws.Range("H40").Value = 2
$ws.Range("I40").Value = 32.299999999999997
$ws.Range("J40").Value = 34.299999999999997

# Row 41: William H Harsha Lake 2016
$ws.Range("A41").Value = "William H Harsha Lake 2016"
$ws.Range("B41").Value = "USEPA"
$ws.Range("C41").Value = 8.3194439599999992
$ws.Range("D41").Value = 34.451219510000001
$ws.Range("E41").Value = 0.095912589
$ws.Range("F41").Value = 11.629
$ws.Range("G41").Value = 0.1552
$ws.Range("H41").Value = 0.6
$ws.Range("I41").Value = 7.7
$ws.Range("J41").Value = 8.3000000000000007

# Rows 42-44: trailing "NA" placeholder rows across columns A:J
foreach ($r in 42..44) {
    foreach ($col in @("A","B","C","D","E","F","G","H","I","J")) {
        $ws.Range("$col$r").Value = "NA"
    }
}

# ---------------------------------------------------------------------------
# View: the sheet is scrolled down to show the newly added rows, and the
# selection moved to C48.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("A2").Select()
$win.FreezePanes = $false
$win.FreezePanes = $true
$win.ScrollRow = 32
$win.ScrollColumn = 1

$ws.Range("C48").Select()
